# Edit script: convert the "Cosmic Voyages / Space-Time" document into the
# "Enchanting World of Biology" document, and normalize the font from the
# (misspelled) "TimesNewToman" to "Times New Roman" throughout.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Fix the font across the whole document (covers every w:rFonts
#    ascii/hAnsi occurrence in one shot).
# ---------------------------------------------------------------------
$docEnd = $d.Content.End
$fullRange = $d.Range(0, $docEnd)
$fullRange.Font.Name = "Times New Roman"

# ---------------------------------------------------------------------
# 2. Title paragraph
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(1)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$r.Text = "The Enchanting World of Biology: Unveiling the Secrets of Life"

# ---------------------------------------------------------------------
# 3. Author name paragraph
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(2)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$r.Text = "Dr. Olivia Green"

# ---------------------------------------------------------------------
# 4. Contact / email paragraph
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(3)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$r.Text = "science.explorer@validmail.com"

# ---------------------------------------------------------------------
# 5. Main body paragraph (paragraph 5 - paragraph 4 is the blank spacer)
# ---------------------------------------------------------------------
$br = [char]11
$bodyText = "Delve into the realm of biology, a fascinating science that unravels the intricate mechanisms of life." +
            " Biology holds the key to understanding the processes that govern living organisms, from the tiniest bacteria to the towering sequoia trees." +
            " In this subject, we embark on a voyage of discovery, exploring the fundamental principles that govern the workings of life, delving into the mysteries of cells, genetics, and ecosystems." +
            $br + $br +
            "Through biology, we uncover the incredible diversity of life on Earth, learning about the astonishing array of organisms that inhabit our planet." +
            " From vibrant coral reefs teeming with marine life to sprawling rainforests teeming with biodiversity, we witness the breathtaking tapestry of life's myriad forms." +
            " With each new discovery, we deepen our understanding of the interconnectedness of all living beings, fostering a sense of wonder and appreciation for the natural world." +
            $br + $br +
            "Moreover, biology holds immense practical relevance, impacting various aspects of human existence." +
            " From developing new medicines to combat diseases, harnessing the potential of genetic engineering to address global challenges, to comprehending the intricate workings of the human body, biology plays a crucial role in improving human health, well-being, and the sustainability of our planet." +
            " By unraveling the intricacies of life, we empower ourselves to make informed decisions, address global challenges, and shape a better future for generations to come."

$p = $d.Paragraphs.Item(5)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$r.Text = $bodyText

# ---------------------------------------------------------------------
# 6. "Summary" heading paragraph - text unchanged, only font (already done)
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 7. Summary body paragraph
# ---------------------------------------------------------------------
$summaryText = "Biology is a captivating science that unveils the secrets of life, guiding us to comprehend the fascinating mechanisms that govern living organisms." +
               " Through the study of cells, genetics, and ecosystems, we unravel the captivating tapestry of diversity that characterizes life on Earth." +
               " Additionally, biology holds immense practical relevance, contributing to advancements in medicine, agriculture, and biotechnology." +
               " By exploring the complexities of life, we acquire valuable insights that empower us to make informed choices, face global challenges, and build a sustainable world for future generations."

$p = $d.Paragraphs.Item(7)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$r.Text = $summaryText

# ---------------------------------------------------------------------
# 8. Add a new empty paragraph at the very end of the document
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRange = $d.Range($p.Range.End, $p.Range.End)
$endRange.InsertParagraphAfter()

Write-Host "Edit complete"
